{"js": "// Replace each two-digit-by-two-digit multiplication problem's text with its\n// new value. Every \"<a>\u00d7<b>=<c>\" string in the table is unique, both before\n// and after the edit, so an exact-text search + replace is unambiguous.\nconst replacements = [\n  [\"29\u00d724=696\", \"63\u00d733=2079\"],\n  [\"98\u00d794=9212\", \"48\u00d729=1392\"],\n  [\"49\u00d756=2744\", \"43\u00d756=2408\"],\n  [\"43\u00d755=2365\", \"50\u00d724=1200\"],\n  [\"80\u00d725=2000\", \"98\u00d724=2352\"],\n  [\"57\u00d758=3306\", \"82\u00d742=3444\"],\n  [\"29\u00d750=1450\", \"79\u00d737=2923\"],\n  [\"86\u00d732=2752\", \"70\u00d774=5180\"],\n  [\"54\u00d721=1134\", \"45\u00d786=3870\"],\n  [\"86\u00d745=3870\", \"96\u00d777=7392\"],\n  [\"13\u00d717=221\", \"85\u00d775=6375\"],\n  [\"56\u00d728=1568\", \"23\u00d775=1725\"],\n  [\"29\u00d778=2262\", \"83\u00d716=1328\"],\n  [\"18\u00d784=1512\", \"40\u00d732=1280\"],\n  [\"70\u00d716=1120\", \"14\u00d747=658\"],\n  [\"45\u00d799=4455\", \"70\u00d744=3080\"],\n  [\"18\u00d745=810\", \"31\u00d741=1271\"],\n  [\"43\u00d778=3354\", \"53\u00d729=1537\"],\n  [\"39\u00d797=3783\", \"89\u00d723=2047\"],\n  [\"42\u00d756=2352\", \"77\u00d717=1309\"],\n  [\"15\u00d713=195\", \"33\u00d735=1155\"],\n  [\"61\u00d746=2806\", \"19\u00d746=874\"],\n  [\"21\u00d796=2016\", \"37\u00d757=2109\"],\n  [\"82\u00d719=1558\", \"26\u00d722=572\"],\n  [\"91\u00d753=4823\", \"17\u00d753=901\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication problem's text with its\n# new value. Every \"<a>\u00d7<b>=<c>\" string in the table is unique, both before\n# and after the edit, so an exact-text Find/Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"29\u00d724=696\", \"63\u00d733=2079\"),\n    @(\"98\u00d794=9212\", \"48\u00d729=1392\"),\n    @(\"49\u00d756=2744\", \"43\u00d756=2408\"),\n    @(\"43\u00d755=2365\", \"50\u00d724=1200\"),\n    @(\"80\u00d725=2000\", \"98\u00d724=2352\"),\n    @(\"57\u00d758=3306\", \"82\u00d742=3444\"),\n    @(\"29\u00d750=1450\", \"79\u00d737=2923\"),\n    @(\"86\u00d732=2752\", \"70\u00d774=5180\"),\n    @(\"54\u00d721=1134\", \"45\u00d786=3870\"),\n    @(\"86\u00d745=3870\", \"96\u00d777=7392\"),\n    @(\"13\u00d717=221\", \"85\u00d775=6375\"),\n    @(\"56\u00d728=1568\", \"23\u00d775=1725\"),\n    @(\"29\u00d778=2262\", \"83\u00d716=1328\"),\n    @(\"18\u00d784=1512\", \"40\u00d732=1280\"),\n    @(\"70\u00d716=1120\", \"14\u00d747=658\"),\n    @(\"45\u00d799=4455\", \"70\u00d744=3080\"),\n    @(\"18\u00d745=810\", \"31\u00d741=1271\"),\n    @(\"43\u00d778=3354\", \"53\u00d729=1537\"),\n    @(\"39\u00d797=3783\", \"89\u00d723=2047\"),\n    @(\"42\u00d756=2352\", \"77\u00d717=1309\"),\n    @(\"15\u00d713=195\", \"33\u00d735=1155\"),\n    @(\"61\u00d746=2806\", \"19\u00d746=874\"),\n    @(\"21\u00d796=2016\", \"37\u00d757=2109\"),\n    @(\"82\u00d719=1558\", \"26\u00d722=572\"),\n    @(\"91\u00d753=4823\", \"17\u00d753=901\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
